# Update row 53's "Time without Visualization" value and append new
# performance-check rows (54-64), matching the latest BDA run data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: row, A (Number of Vertex), B (Time without Visualization)
$data = @(
    @(53, 8, "1.157407407407407e-08"),
    @(54, 8, "3.472222222222222e-08"),
    @(55, 8, "1.157407407407407e-08"),
    @(56, 8, "2.314814814814815e-08"),
    @(57, 8, "2.314814814814815e-08"),
    @(58, 7, "1.157407407407407e-08"),
    @(59, 7, "7.638888888888889e-07"),
    @(60, 8, "2.314814814814815e-08"),
    @(61, 8, "2.314814814814815e-08"),
    @(62, 8, "1.157407407407407e-08"),
    @(63, 8, "8.564814814814814e-07"),
    @(64, 8, "2.297453703703704e-08")
)

foreach ($item in $data) {
    $rowNum = $item[0]
    $vertexCount = $item[1]
    $timeValue = [double]$item[2]

    $ws.Cells.Item($rowNum, 1).Value = $vertexCount
    $ws.Cells.Item($rowNum, 2).Value = $timeValue
    $ws.Cells.Item($rowNum, 2).NumberFormat = "[hh]:mm:ss"
    $ws.Cells.Item($rowNum, 3).Value = "-"
}
